$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 3 onto row 4 so the new entry gets the same
# per-column cell styles (text style, date style, date style, time style)
# without Excel having to fabricate brand-new merged styles.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)

# Fill in the new bitácora entry.
$ws.Range("A4").Value = "Se realizó un cambio en la interfaz de clientes y se creó la interfaz de gestión de cuentas."
$ws.Range("B4").Value = 42864.916666666664
$ws.Range("C4").Value = 42865.104166666664
$ws.Range("D4").Value = 0.1875

# Row 4 needs the same (taller) row height as the other populated rows so the
# wrapped description text is fully visible.
$ws.Rows.Item(4).RowHeight = 30

# Update the current selection to match the author's final cursor position.
$ws.Range("B5").Select()
